# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '27.060.34'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  +0.50%  '
$ws.Cells.Item(2, 5).Style = 'Normal'
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.675.01'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  +0.23%  '
$ws.Cells.Item(3, 5).Style = 'Normal'
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  +0.07%  '
$ws.Cells.Item(4, 5).Style = 'Normal'
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '215.04'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +0.15%  '
$ws.Cells.Item(5, 5).Style = 'Normal'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  -0.22%  '
$ws.Cells.Item(6, 5).Style = 'Normal'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$ws.Cells.Item(7, 5).Style = 'Normal'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.255'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  +1.93%  '
$ws.Cells.Item(8, 5).Style = 'Normal'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  +0.09%  '
$ws.Cells.Item(9, 5).Style = 'Normal'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '21.20'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +4.77%  '
$ws.Cells.Item(10, 5).Style = 'Normal'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0883'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  -0.80%  '
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '1.910.65'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.665.78'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  -0.27%  '
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '4.12'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  +0.91%  '
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  +1.53%  '
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '65.99'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  +0.77%  '
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '27.032.91'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  +0.39%  '
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(18, 2).Value = 'BitcoinCash'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '237.28'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  +1.74%  '
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '8.15'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  +1.61%  '
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  +0.96%  '
$ws.Cells.Item(20, 5).Style = 'Normal'
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '4.46'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '9.34'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  +2.23%  '
$ws.Cells.Item(23, 5).Style = 'Normal'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  -1.95%  '
$ws.Cells.Item(24, 5).Style = 'Normal'
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '146.19'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  -0.06%  '
$ws.Cells.Item(25, 5).Style = 'Normal'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  +1.59%  '
$ws.Cells.Item(26, 5).Style = 'Normal'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '16.35'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  +2.71%  '
$ws.Cells.Item(27, 5).Style = 'Normal'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.112'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  +0.51%  '
$ws.Cells.Item(28, 5).Style = 'Normal'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  -0.08%  '
$ws.Cells.Item(29, 5).Style = 'Normal'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.0498'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  +0.01%  '
$ws.Cells.Item(30, 5).Style = 'Normal'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -0.05%  '
$ws.Cells.Item(31, 5).Style = 'Normal'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.35'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  +0.96%  '
$ws.Cells.Item(32, 5).Style = 'Normal'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.541.39'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +6.15%  '
$ws.Cells.Item(33, 5).Style = 'Normal'
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  +1.90%  '
$ws.Cells.Item(34, 5).Style = 'Normal'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  +3.78%  '
$ws.Cells.Item(35, 5).Style = 'Normal'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -1.02%  '
$ws.Cells.Item(36, 5).Style = 'Normal'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.597'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  +3.07%  '
$ws.Cells.Item(37, 5).Style = 'Normal'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.921'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  +2.34%  '
$ws.Cells.Item(38, 5).Style = 'Normal'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  +2.21%  '
$ws.Cells.Item(39, 5).Style = 'Normal'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  +2.26%  '
$ws.Cells.Item(40, 5).Style = 'Normal'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  +0.00%  '
$ws.Cells.Item(41, 5).Style = 'Normal'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '67.68'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  +2.19%  '
$ws.Cells.Item(42, 5).Style = 'Normal'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '5.58'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -2.80%  '
$ws.Cells.Item(43, 5).Style = 'Normal'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.26'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -1.75%  '
$ws.Cells.Item(44, 5).Style = 'Normal'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '1.818.53'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  +0.67%  '
$ws.Cells.Item(45, 5).Style = 'Normal'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.785'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  +0.59%  '
$ws.Cells.Item(46, 5).Style = 'Normal'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '90.83'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  +0.23%  '
$ws.Cells.Item(47, 5).Style = 'Normal'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  +1.82%  '
$ws.Cells.Item(48, 5).Style = 'Normal'
$ws.Cells.Item(49, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.0₆0106'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  +1.06%  '
$ws.Cells.Item(49, 5).Style = 'Normal'
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.104'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  +2.53%  '
$ws.Cells.Item(50, 5).Style = 'Normal'
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '8.03'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +5.48%  '
$ws.Cells.Item(51, 5).Style = 'Normal'
